$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# Schedule sheet updates
$ws1.Range("E2").Value = 768.514968
$ws1.Range("F2").Value = 12.70692738095238

# Detailed sheet updates
$ws2.Range("B2").Value = 57.09
$ws2.Range("B3").Value = 57.09
$ws2.Range("B4").Value = 57.09
$ws2.Range("C4").Value = "historical"
$ws2.Range("B6").Value = 57.06003
$ws2.Range("B7").Value = 40.54
$ws2.Range("B8").Value = 40.54
$ws2.Range("B9").Value = 56.98
$ws2.Range("B10").Value = 56.98
$ws2.Range("B14").Value = 57.06018
$ws2.Range("B15").Value = 50.75171
$ws2.Range("B16").Value = 36.2
$ws2.Range("B17").Value = 36.06029
$ws2.Range("B18").Value = 15.43474
$ws2.Range("B19").Value = 13.96726
$ws2.Range("B20").Value = 36.06011
$ws2.Range("B21").Value = 35.88
$ws2.Range("B22").Value = 36.06057
$ws2.Range("B24").Value = 36.06092
$ws2.Range("B25").Value = 36.06046
$ws2.Range("B27").Value = 36.06
$ws2.Range("B28").Value = 36.0609
$ws2.Range("B29").Value = 31.0352
$ws2.Range("B30").Value = 36.06029
$ws2.Range("B31").Value = 36.06031
$ws2.Range("B35").Value = 22.50263
$ws2.Range("B36").Value = -0.57355
$ws2.Range("B37").Value = -2.99308
$ws2.Range("B38").Value = -2.83044
$ws2.Range("B39").Value = -2.77762
$ws2.Range("B41").Value = 8.580550000000001
$ws2.Range("B42").Value = 11.92003
$ws2.Range("B43").Value = 11.09544
$ws2.Range("B44").Value = 9.361459999999999
$ws2.Range("B45").Value = 9.751200000000001
$ws2.Range("B47").Value = 56.98
$ws2.Range("B49").Value = 56.98
